$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '28.451.10'
Set-TextValue $ws.Range("E2") '  +3.59%  '
Set-TextValue $ws.Range("D3") '1.869.97'
Set-TextValue $ws.Range("E3") '  +2.00%  '
Set-TextValue $ws.Range("E4") '  -0.13%  '
Set-TextValue $ws.Range("D5") '338.86'
Set-TextValue $ws.Range("E5") '  +2.24%  '
Set-TextValue $ws.Range("D6") '1.000'
Set-TextValue $ws.Range("E6") '  -0.13%  '
Set-TextValue $ws.Range("D7") '0.4708'
Set-TextValue $ws.Range("E7") '  +2.18%  '
Set-TextValue $ws.Range("D8") '0.3976'
Set-TextValue $ws.Range("E8") '  +3.79%  '
Set-TextValue $ws.Range("D9") '47.72'
Set-TextValue $ws.Range("E9") '  +2.43%  '
Set-TextValue $ws.Range("D10") '0.08027'
Set-TextValue $ws.Range("E10") '  +1.56%  '
Set-TextValue $ws.Range("E11") '  +2.97%  '
Set-TextValue $ws.Range("D12") '22.10'
Set-TextValue $ws.Range("E12") '  +4.81%  '
Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '6.035'
Set-TextValue $ws.Range("E13") '  +2.51%  '
Set-TextValue $ws.Range("B14") 'WrappedEther'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D14") '1.874.26'
Set-TextValue $ws.Range("E14") '  +3.25%  '
Set-TextValue $ws.Range("D15") '7.290'
Set-TextValue $ws.Range("E15") '  +3.32%  '
Set-TextValue $ws.Range("D16") '90.92'
Set-TextValue $ws.Range("E16") '  +3.28%  '
Set-TextValue $ws.Range("E17") '  -0.11%  '
Set-TextValue $ws.Range("D18") '0.00001044'
Set-TextValue $ws.Range("E18") '  +1.50%  '
Set-TextValue $ws.Range("D19") '0.06622'
Set-TextValue $ws.Range("E19") '  +0.01%  '
Set-TextValue $ws.Range("D20") '17.57'
Set-TextValue $ws.Range("E20") '  +3.01%  '
Set-TextValue $ws.Range("D21") '1.000'
Set-TextValue $ws.Range("E21") '  -0.18%  '
Set-TextValue $ws.Range("D22") '28.453.63'
Set-TextValue $ws.Range("E22") '  +3.58%  '
Set-TextValue $ws.Range("D23") '5.474'
Set-TextValue $ws.Range("E23") '  +2.40%  '
Set-TextValue $ws.Range("D24") '11.06'
Set-TextValue $ws.Range("E24") '  +2.15%  '
Set-TextValue $ws.Range("D25") '2.268'
Set-TextValue $ws.Range("E25") '  -1.88%  '
Set-TextValue $ws.Range("D26") '2.091.62'
Set-TextValue $ws.Range("E26") '  +2.73%  '
Set-TextValue $ws.Range("D27") '160.78'
Set-TextValue $ws.Range("E27") '  +2.16%  '
Set-TextValue $ws.Range("D28") '19.82'
Set-TextValue $ws.Range("E28") '  +2.13%  '
Set-TextValue $ws.Range("D29") '2.121'
Set-TextValue $ws.Range("E29") '  +2.66%  '
Set-TextValue $ws.Range("D30") '5.510'
Set-TextValue $ws.Range("E30") '  +3.96%  '
Set-TextValue $ws.Range("D31") '120.31'
Set-TextValue $ws.Range("E31") '  +1.00%  '
Set-TextValue $ws.Range("D32") '0.9793'
Set-TextValue $ws.Range("E32") '  +2.45%  '
Set-TextValue $ws.Range("D33") '0.09514'
Set-TextValue $ws.Range("E33") '  +2.51%  '
Set-TextValue $ws.Range("D34") '3.586'
Set-TextValue $ws.Range("E34") '  +0.24%  '
Set-TextValue $ws.Range("D35") '1.377'
Set-TextValue $ws.Range("E35") '  +4.79%  '
Set-TextValue $ws.Range("D36") '5.362'
Set-TextValue $ws.Range("E36") '  +2.37%  '
Set-TextValue $ws.Range("D37") '0.06117'
Set-TextValue $ws.Range("E37") '  +3.04%  '
Set-TextValue $ws.Range("D38") '0.02258'
Set-TextValue $ws.Range("E38") '  +3.11%  '
Set-TextValue $ws.Range("D39") '8.369'
Set-TextValue $ws.Range("E39") '  +3.64%  '
Set-TextValue $ws.Range("D40") '1.183'
Set-TextValue $ws.Range("E40") '  +1.97%  '
Set-TextValue $ws.Range("D41") '0.5949'
Set-TextValue $ws.Range("E41") '  +2.58%  '
Set-TextValue $ws.Range("D42") '0.9998'
Set-TextValue $ws.Range("E42") '  -0.10%  '
Set-TextValue $ws.Range("D43") '0.1879'
Set-TextValue $ws.Range("E43") '  +2.15%  '
Set-TextValue $ws.Range("E44") '  +3.30%  '
Set-TextValue $ws.Range("D45") '1.283'
Set-TextValue $ws.Range("E45") '  +0.10%  '
Set-TextValue $ws.Range("D46") '0.5587'
Set-TextValue $ws.Range("E46") '  +1.83%  '
Set-TextValue $ws.Range("D47") '12.20'
Set-TextValue $ws.Range("E47") '  +1.61%  '
Set-TextValue $ws.Range("D48") '1.959'
Set-TextValue $ws.Range("E48") '  +4.67%  '
Set-TextValue $ws.Range("D49") '0.06968'
Set-TextValue $ws.Range("E49") '  +4.78%  '
Set-TextValue $ws.Range("D50") '2.081'
Set-TextValue $ws.Range("E50") '  +15.58%  '
Set-TextValue $ws.Range("D51") '111.86'
Set-TextValue $ws.Range("E51") '  +1.29%  '
